# Backlog.xlsx - "Fix ingreso de materiales por orden de compra"
#
# 1) Mark a handful of already-finished backlog items as "terminado"
#    (they were still showing "no comenzado").
# 2) Log two new backlog items ("no comenzado") at the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Hoja1" (Backlog) is the active/tabSelected sheet

# --- Mark tasks as finished ---------------------------------------------
$doneRows = 67, 87, 88, 90, 93
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 2).Value = "terminado"
}

# --- Append new backlog entries -----------------------------------------
$newTasks = @(
    "facturacion, permitir modificar liquidacion de iibb manualmente",
    "reporte para saber en que piezas es utilizado un componente"
)

$lastRow = 98
for ($i = 0; $i -lt $newTasks.Count; $i++) {
    $row = $lastRow + 1 + $i
    $ws.Cells.Item($row, 1).Value = $newTasks[$i]
    $ws.Cells.Item($row, 2).Value = "no comenzado"
}

# --- Match the author's on-screen selection when they saved -------------
$ws.Range("C97").Select()
